# Weekly fruit/vegetable price update:
# Insert a new daily record as row 56 (pushing the existing rows 56-87 down
# to 57-88) and populate it with the new week's figures while keeping the
# other attributes identical to the record that used to sit in row 56
# (which is now row 57).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 56; everything below (old rows
# 56-87) shifts down to 57-88, carrying its own formatting/values with it.
$ws.Rows.Item(56).Insert()

# Duplicate the record that is now in row 57 into the new blank row 56, so
# row 56 starts out identical to it ...
$ws.Range("A57:R57").Copy()
$ws.Range("A56:R56").PasteSpecial()

# ... then overwrite this new record's date (D) and volume (J) with the
# values for the new week.
$ws.Range("D56").Value = 44824
$ws.Range("J56").Value = 20
